$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear old data body (rows 2-25), keep header row 1
$ws.Range("A2:H25").ClearContents()

# Write full refreshed results table (rows 2-28)
$ws.Cells.Item(2, 1).Value = "ridge_regression"
$ws.Cells.Item(2, 2).Value = "rfe"
$ws.Cells.Item(2, 3).Value = "{'alpha': 1.0, 'fit_intercept': True, 'max_iter': 100}"
$ws.Cells.Item(2, 4).Value = 139.976860820079
$ws.Cells.Item(2, 5).Value = 7.12835557099886
$ws.Cells.Item(2, 6).Value = 0.420500265262409
$ws.Cells.Item(2, 7).Value = 10.1373035611276
$ws.Cells.Item(2, 8).Value = 0.433322517241577

$ws.Cells.Item(3, 1).Value = "ridge_regression"
$ws.Cells.Item(3, 2).Value = "mutual information"
$ws.Cells.Item(3, 3).Value = "{'alpha': 1.0, 'fit_intercept': True, 'max_iter': 100}"
$ws.Cells.Item(3, 4).Value = 142.228576521385
$ws.Cells.Item(3, 5).Value = 7.17183476313198
$ws.Cells.Item(3, 6).Value = 0.434481539959661
$ws.Cells.Item(3, 7).Value = 10.100022640484
$ws.Cells.Item(3, 8).Value = 0.451733229199198

$ws.Cells.Item(4, 1).Value = "ridge_regression"
$ws.Cells.Item(4, 2).Value = "random forest"
$ws.Cells.Item(4, 3).Value = "{'alpha': 1.0, 'fit_intercept': True, 'max_iter': 100}"
$ws.Cells.Item(4, 4).Value = 141.284802060638
$ws.Cells.Item(4, 5).Value = 7.16453480514923
$ws.Cells.Item(4, 6).Value = 0.434740973861546
$ws.Cells.Item(4, 7).Value = 10.0888496598485
$ws.Cells.Item(4, 8).Value = 0.452019604292011

$ws.Cells.Item(5, 1).Value = "ridge_regression"
$ws.Cells.Item(5, 2).Value = "all"
$ws.Cells.Item(5, 3).Value = "{'alpha': 10.0, 'fit_intercept': True, 'max_iter': 100}"
$ws.Cells.Item(5, 4).Value = 136.873793984982
$ws.Cells.Item(5, 5).Value = 7.05638092001734
$ws.Cells.Item(5, 6).Value = 0.451442322824882
$ws.Cells.Item(5, 7).Value = 9.93564557424806
$ws.Cells.Item(5, 8).Value = 0.468384856771993

$ws.Cells.Item(6, 1).Value = "ridge_regression"
$ws.Cells.Item(6, 2).Value = "rfe_corr"
$ws.Cells.Item(6, 3).Value = "{'alpha': 10.0, 'fit_intercept': True, 'max_iter': 100}"
$ws.Cells.Item(6, 4).Value = 126.712044726729
$ws.Cells.Item(6, 5).Value = 7.03175815779687
$ws.Cells.Item(6, 6).Value = 0.476649025726928
$ws.Cells.Item(6, 7).Value = 9.65105816857169
$ws.Cells.Item(6, 8).Value = 0.494627366823451

$ws.Cells.Item(7, 1).Value = "ridge_regression"
$ws.Cells.Item(7, 2).Value = "mutual information_corr"
$ws.Cells.Item(7, 3).Value = "{'alpha': 10.0, 'fit_intercept': True, 'max_iter': 100}"
$ws.Cells.Item(7, 4).Value = 125.240102033029
$ws.Cells.Item(7, 5).Value = 6.89591334259609
$ws.Cells.Item(7, 6).Value = 0.488749044494408
$ws.Cells.Item(7, 7).Value = 9.56623183250049
$ws.Cells.Item(7, 8).Value = 0.508888712547718

$ws.Cells.Item(8, 1).Value = "ridge_regression"
$ws.Cells.Item(8, 2).Value = "random forest_corr"
$ws.Cells.Item(8, 3).Value = "{'alpha': 10.0, 'fit_intercept': True, 'max_iter': 100}"
$ws.Cells.Item(8, 4).Value = 123.628126479426
$ws.Cells.Item(8, 5).Value = 6.87151821858771
$ws.Cells.Item(8, 6).Value = 0.490440098703249
$ws.Cells.Item(8, 7).Value = 9.53179573550928
$ws.Cells.Item(8, 8).Value = 0.510559818735047

$ws.Cells.Item(9, 1).Value = "ridge_regression"
$ws.Cells.Item(9, 2).Value = "all_corr"
$ws.Cells.Item(9, 3).Value = "{'alpha': 10.0, 'fit_intercept': True, 'max_iter': 100}"
$ws.Cells.Item(9, 4).Value = 123.429318443384
$ws.Cells.Item(9, 5).Value = 6.84016828261159
$ws.Cells.Item(9, 6).Value = 0.492165622300035
$ws.Cells.Item(9, 7).Value = 9.51941929893485
$ws.Cells.Item(9, 8).Value = 0.511472669309067

$ws.Cells.Item(10, 1).Value = "ridge_regression"
$ws.Cells.Item(10, 2).Value = "authors"
$ws.Cells.Item(10, 3).Value = "{'alpha': 1.0, 'fit_intercept': True, 'max_iter': 100}"
$ws.Cells.Item(10, 4).Value = 139.979973592673
$ws.Cells.Item(10, 5).Value = 7.24694579931723
$ws.Cells.Item(10, 6).Value = 0.431266965605833
$ws.Cells.Item(10, 7).Value = 10.0841257489019
$ws.Cells.Item(10, 8).Value = 0.447361187582569

$ws.Cells.Item(11, 1).Value = "random_forest"
$ws.Cells.Item(11, 2).Value = "rfe"
$ws.Cells.Item(11, 3).Value = "{'max_depth': 8, 'max_features': 'sqrt', 'n_estimators': 200}"
$ws.Cells.Item(11, 4).Value = 81.3031889456529
$ws.Cells.Item(11, 5).Value = 5.35936682471787
$ws.Cells.Item(11, 6).Value = 0.60525287143655
$ws.Cells.Item(11, 7).Value = 8.01306746203909
$ws.Cells.Item(11, 8).Value = 0.607406118711081

$ws.Cells.Item(12, 1).Value = "random_forest"
$ws.Cells.Item(12, 2).Value = "mutual information"
$ws.Cells.Item(12, 3).Value = "{'max_depth': 5, 'max_features': 'sqrt', 'n_estimators': 200}"
$ws.Cells.Item(12, 4).Value = 78.1342926689553
$ws.Cells.Item(12, 5).Value = 5.57042843237904
$ws.Cells.Item(12, 6).Value = 0.618602989880979
$ws.Cells.Item(12, 7).Value = 7.89094524270069
$ws.Cells.Item(12, 8).Value = 0.623838686207411

$ws.Cells.Item(13, 1).Value = "random_forest"
$ws.Cells.Item(13, 2).Value = "random forest"
$ws.Cells.Item(13, 3).Value = "{'max_depth': 5, 'max_features': 'sqrt', 'n_estimators': 200}"
$ws.Cells.Item(13, 4).Value = 78.4149526730382
$ws.Cells.Item(13, 5).Value = 5.61437480589397
$ws.Cells.Item(13, 6).Value = 0.62011383850581
$ws.Cells.Item(13, 7).Value = 7.89516540178532
$ws.Cells.Item(13, 8).Value = 0.625012738085308

$ws.Cells.Item(14, 1).Value = "random_forest"
$ws.Cells.Item(14, 2).Value = "all"
$ws.Cells.Item(14, 3).Value = "{'max_depth': 5, 'max_features': 'sqrt', 'n_estimators': 200}"
$ws.Cells.Item(14, 4).Value = 80.857407321346
$ws.Cells.Item(14, 5).Value = 5.7268683029627
$ws.Cells.Item(14, 6).Value = 0.614058145740051
$ws.Cells.Item(14, 7).Value = 7.99086294545401
$ws.Cells.Item(14, 8).Value = 0.619221037458723

$ws.Cells.Item(15, 1).Value = "random_forest"
$ws.Cells.Item(15, 2).Value = "rfe_corr"
$ws.Cells.Item(15, 3).Value = "{'max_depth': 5, 'max_features': 'log2', 'n_estimators': 200}"
$ws.Cells.Item(15, 4).Value = 80.575123706367
$ws.Cells.Item(15, 5).Value = 5.70053979876665
$ws.Cells.Item(15, 6).Value = 0.606092959782884
$ws.Cells.Item(15, 7).Value = 8.01231372752469
$ws.Cells.Item(15, 8).Value = 0.610907860266208

$ws.Cells.Item(16, 1).Value = "random_forest"
$ws.Cells.Item(16, 2).Value = "mutual information_corr"
$ws.Cells.Item(16, 3).Value = "{'max_depth': 5, 'max_features': 'sqrt', 'n_estimators': 200}"
$ws.Cells.Item(16, 4).Value = 79.9615162062374
$ws.Cells.Item(16, 5).Value = 5.69962878896615
$ws.Cells.Item(16, 6).Value = 0.612320435753989
$ws.Cells.Item(16, 7).Value = 7.96852160114825
$ws.Cells.Item(16, 8).Value = 0.616898403028324

$ws.Cells.Item(17, 1).Value = "random_forest"
$ws.Cells.Item(17, 2).Value = "random forest_corr"
$ws.Cells.Item(17, 3).Value = "{'max_depth': 5, 'max_features': 'sqrt', 'n_estimators': 200}"
$ws.Cells.Item(17, 4).Value = 80.4759512703862
$ws.Cells.Item(17, 5).Value = 5.74730512407592
$ws.Cells.Item(17, 6).Value = 0.613737098294629
$ws.Cells.Item(17, 7).Value = 7.97538864404933
$ws.Cells.Item(17, 8).Value = 0.618733403498939

$ws.Cells.Item(18, 1).Value = "random_forest"
$ws.Cells.Item(18, 2).Value = "all_corr"
$ws.Cells.Item(18, 3).Value = "{'max_depth': 5, 'max_features': 'sqrt', 'n_estimators': 200}"
$ws.Cells.Item(18, 4).Value = 79.9134999770499
$ws.Cells.Item(18, 5).Value = 5.65861359845766
$ws.Cells.Item(18, 6).Value = 0.6137527241857
$ws.Cells.Item(18, 7).Value = 7.96170769298696
$ws.Cells.Item(18, 8).Value = 0.618229597065856

$ws.Cells.Item(19, 1).Value = "random_forest"
$ws.Cells.Item(19, 2).Value = "authors"
$ws.Cells.Item(19, 3).Value = "{'max_depth': 5, 'max_features': 'sqrt', 'n_estimators': 200}"
$ws.Cells.Item(19, 4).Value = 84.0584879625577
$ws.Cells.Item(19, 5).Value = 5.90058585763032
$ws.Cells.Item(19, 6).Value = 0.601543134372753
$ws.Cells.Item(19, 7).Value = 8.12971140281904
$ws.Cells.Item(19, 8).Value = 0.606948528980529

$ws.Cells.Item(20, 1).Value = "neural_network"
$ws.Cells.Item(20, 2).Value = "rfe"
$ws.Cells.Item(20, 3).Value = "{'activation': 'relu', 'hidden_layer_sizes': [128, 64], 'learning_rate': 'adaptive', 'max_iter': 5000, 'solver': 'sgd'}"
$ws.Cells.Item(20, 4).Value = 79.2493132684659
$ws.Cells.Item(20, 5).Value = 5.30246646814824
$ws.Cells.Item(20, 6).Value = 0.618443423794228
$ws.Cells.Item(20, 7).Value = 7.92450618948622
$ws.Cells.Item(20, 8).Value = 0.62181993673031

$ws.Cells.Item(21, 1).Value = "neural_network"
$ws.Cells.Item(21, 2).Value = "mutual information"
$ws.Cells.Item(21, 3).Value = "{'activation': 'relu', 'hidden_layer_sizes': [128, 64], 'learning_rate': 'adaptive', 'max_iter': 5000, 'solver': 'sgd'}"
$ws.Cells.Item(21, 4).Value = 83.0532339173407
$ws.Cells.Item(21, 5).Value = 5.37730930544781
$ws.Cells.Item(21, 6).Value = 0.616310405295683
$ws.Cells.Item(21, 7).Value = 8.02426474556003
$ws.Cells.Item(21, 8).Value = 0.620480313004877

$ws.Cells.Item(22, 1).Value = "neural_network"
$ws.Cells.Item(22, 2).Value = "random forest"
$ws.Cells.Item(22, 3).Value = "{'activation': 'relu', 'hidden_layer_sizes': [20], 'learning_rate': 'adaptive', 'max_iter': 5000, 'solver': 'sgd'}"
$ws.Cells.Item(22, 4).Value = 81.5631873319203
$ws.Cells.Item(22, 5).Value = 5.38466607185225
$ws.Cells.Item(22, 6).Value = 0.61852873806166
$ws.Cells.Item(22, 7).Value = 7.98555799196988
$ws.Cells.Item(22, 8).Value = 0.6223449920808

$ws.Cells.Item(23, 1).Value = "neural_network"
$ws.Cells.Item(23, 2).Value = "all"
$ws.Cells.Item(23, 3).Value = "{'activation': 'relu', 'hidden_layer_sizes': [128, 64], 'learning_rate': 'adaptive', 'max_iter': 5000, 'solver': 'sgd'}"
$ws.Cells.Item(23, 4).Value = 78.0835715702577
$ws.Cells.Item(23, 5).Value = 5.23506725788577
$ws.Cells.Item(23, 6).Value = 0.635634087223037
$ws.Cells.Item(23, 7).Value = 7.80283561667466
$ws.Cells.Item(23, 8).Value = 0.639455020801961

$ws.Cells.Item(24, 1).Value = "neural_network"
$ws.Cells.Item(24, 2).Value = "rfe_corr"
$ws.Cells.Item(24, 3).Value = "{'activation': 'relu', 'hidden_layer_sizes': [20], 'learning_rate': 'constant', 'max_iter': 5000, 'solver': 'adam'}"
$ws.Cells.Item(24, 4).Value = 92.6385076861146
$ws.Cells.Item(24, 5).Value = 6.13258507669207
$ws.Cells.Item(24, 6).Value = 0.560734271513309
$ws.Cells.Item(24, 7).Value = 8.55007096172774
$ws.Cells.Item(24, 8).Value = 0.570504650038458

$ws.Cells.Item(25, 1).Value = "neural_network"
$ws.Cells.Item(25, 2).Value = "mutual information_corr"
$ws.Cells.Item(25, 3).Value = "{'activation': 'relu', 'hidden_layer_sizes': [128, 64], 'learning_rate': 'adaptive', 'max_iter': 5000, 'solver': 'sgd'}"
$ws.Cells.Item(25, 4).Value = 82.0062100260137
$ws.Cells.Item(25, 5).Value = 5.36977308508505
$ws.Cells.Item(25, 6).Value = 0.619685627425278
$ws.Cells.Item(25, 7).Value = 7.98781792060921
$ws.Cells.Item(25, 8).Value = 0.622545765923103

$ws.Cells.Item(26, 1).Value = "neural_network"
$ws.Cells.Item(26, 2).Value = "random forest_corr"
$ws.Cells.Item(26, 3).Value = "{'activation': 'relu', 'hidden_layer_sizes': [20], 'learning_rate': 'constant', 'max_iter': 5000, 'solver': 'adam'}"
$ws.Cells.Item(26, 4).Value = 91.8519165850961
$ws.Cells.Item(26, 5).Value = 5.99820768412173
$ws.Cells.Item(26, 6).Value = 0.579352525502742
$ws.Cells.Item(26, 7).Value = 8.43418736266557
$ws.Cells.Item(26, 8).Value = 0.584292939042452

$ws.Cells.Item(27, 1).Value = "neural_network"
$ws.Cells.Item(27, 2).Value = "all_corr"
$ws.Cells.Item(27, 3).Value = "{'activation': 'relu', 'hidden_layer_sizes': [128, 64], 'learning_rate': 'adaptive', 'max_iter': 5000, 'solver': 'sgd'}"
$ws.Cells.Item(27, 4).Value = 82.215822878305
$ws.Cells.Item(27, 5).Value = 5.36265804249786
$ws.Cells.Item(27, 6).Value = 0.622171492840462
$ws.Cells.Item(27, 7).Value = 7.97926927200173
$ws.Cells.Item(27, 8).Value = 0.624761648242096

$ws.Cells.Item(28, 1).Value = "neural_network"
$ws.Cells.Item(28, 2).Value = "authors"
$ws.Cells.Item(28, 3).Value = "{'activation': 'relu', 'hidden_layer_sizes': [50, 25], 'learning_rate': 'adaptive', 'max_iter': 5000, 'solver': 'sgd'}"
$ws.Cells.Item(28, 4).Value = 85.5774535747547
$ws.Cells.Item(28, 5).Value = 5.58367637073834
$ws.Cells.Item(28, 6).Value = 0.590909579370096
$ws.Cells.Item(28, 7).Value = 8.22410861103194
$ws.Cells.Item(28, 8).Value = 0.597199347212664

# Column C width update
$ws.Columns("C:C").ColumnWidth = 90.97

# Reset selection to A1
$ws.Range("A1").Select()